{"js": "// The edit merges each of the Title / Author / Abstract paragraphs\n// (previously split into one run per word/space) into a single run\n// containing the full paragraph text, with no other content changes.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style\");\nawait context.sync();\n\nconst replacements = {\n  \"Title\": \"Questions: Introduction to simultaneous equations\",\n  \"Author\": \"Ollie Brooke\",\n  \"Abstract\": \"Questions relating to the introduction to simultaneous equations study guide.\"\n};\n\nfor (const para of paragraphs.items) {\n  if (Object.prototype.hasOwnProperty.call(replacements, para.style)) {\n    para.insertText(replacements[para.style], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# The edit merges each of the Title / Author / Abstract paragraphs\n# (previously split into one run per word/space) into a single run\n# containing the full paragraph text, with no other content changes.\n$d = $word.ActiveDocument\n\n$replacements = @{\n  \"Title\"    = \"Questions: Introduction to simultaneous equations\"\n  \"Author\"   = \"Ollie Brooke\"\n  \"Abstract\" = \"Questions relating to the introduction to simultaneous equations study guide.\"\n}\n\nforeach ($p in $d.Paragraphs) {\n  $styleName = $p.Style.NameLocal\n  if ($replacements.ContainsKey($styleName)) {\n    $newText = $replacements[$styleName]\n    $r = $p.Range\n    $r.MoveEnd(1, -1) | Out-Null\n    $r.Find.Execute($r.Text, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n  }\n}\n"}
